# Added implementation of MSM measure.
# The "interfaceOperations" sheet previously only listed the operations
# explicitly declared on each interface. It now also lists the operations
# inherited from java.lang.Object that the interface's (only) implementation
# surfaces, matching the richer listing already used on "classMethods".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("interfaceOperations")

$interfaceName = "com.zatribune.spring.ecommerce.orders.controller.OrderController"

$rows = @(
    @($interfaceName, "equals(java.lang.Object)", "public", "boolean"),
    @($interfaceName, "toString()", "public", "java.lang.String"),
    @($interfaceName, "all()", "public", "java.util.List"),
    @($interfaceName, "getClass()", "public", "java.lang.Class"),
    @($interfaceName, "OrderController(org.springframework.kafka.core.KafkaTemplate, org.springframework.kafka.config.StreamsBuilderFactoryBean)", "public", "void"),
    @($interfaceName, "notifyAll()", "public", "void"),
    @($interfaceName, "hashCode()", "public", "int"),
    @($interfaceName, "wait()", "public", "void"),
    @($interfaceName, "notify()", "public", "void"),
    @($interfaceName, "create(domain.Order)", "public", "domain.Order"),
    @($interfaceName, "wait(long)", "public", "void"),
    @($interfaceName, "wait(long, int)", "public", "void")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
